$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Summary" - update aggregate metrics after trade #89 closed
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("B3").Value = 1199.74    # Current Capital
$summary.Range("B4").Value = -0.25      # Total P&L $
$summary.Range("B5").Value = -0.06      # Total P&L %
$summary.Range("B6").Value = 89         # Total Trades
$summary.Range("B8").Value = 35         # Losing Trades
$summary.Range("B9").Value = 42.7       # Win Rate %

# ---------------------------------------------------------------------
# Sheet 2: "Strategy Status" - update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item(2)
$status.Range("C4").Value = 99.73999999999999   # Capital
$status.Range("D4").Value = 89                  # Trades
$status.Range("E4").Value = -0.25               # P&L $
$status.Range("F4").Value = -0.26               # P&L %
$status.Range("G4").Value = 42.7                # Win Rate %

# ---------------------------------------------------------------------
# Helper: append the new trade #89 row to a trades sheet at row 90
# ---------------------------------------------------------------------
function Add-Trade89Row($ws) {
    $ws.Range("A90").Value = 89

    # Date/time columns must remain plain text (not auto-converted to
    # date/time serials) to match the original inline-string storage.
    $ws.Range("B90").NumberFormat = "@"
    $ws.Range("B90").Value = "2026-02-17"
    $ws.Range("C90").Value = "09:08:21"

    $ws.Range("D90").Value = "MarketMaking"
    $ws.Range("E90").Value = "DOWN"
    $ws.Range("F90").Value = 0.31
    $ws.Range("G90").Value = 0.18
    $ws.Range("H90").Value = "CLOSED"
    $ws.Range("I90").Value = -41.9355
    $ws.Range("J90").Value = -0.13
    $ws.Range("K90").Value = 99.73999999999999
    $ws.Range("L90").Value = 0
    $ws.Range("M90").Value = 0
    $ws.Range("N90").Value = 0.6
    $ws.Range("O90").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P90").Value = "early_exit"
    $ws.Range("Q90").Value = 0.14
}

# ---------------------------------------------------------------------
# Sheet 3: "All Trades" - append new trade row
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item(3)
Add-Trade89Row $allTrades

# ---------------------------------------------------------------------
# Sheet 4: "MarketMaking" - append the same new trade row
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item(4)
Add-Trade89Row $marketMaking
